$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 94868
$ws.Range("C2").Value = 121851
$ws.Range("D2").Value = 160358
$ws.Range("E2").Value = 271853
$ws.Range("F2").Value = 469288
$ws.Range("G2").Value = 82624
$ws.Range("H2").Value = 71270
$ws.Range("B3").Value = 95281
$ws.Range("C3").Value = 131000
$ws.Range("D3").Value = 161876
$ws.Range("E3").Value = 279336
$ws.Range("F3").Value = 470741
$ws.Range("G3").Value = 82841
$ws.Range("H3").Value = 71550
$ws.Range("B4").Value = 107182
$ws.Range("C4").Value = 131268
$ws.Range("D4").Value = 182403
$ws.Range("E4").Value = 314411
$ws.Range("F4").Value = 530037
$ws.Range("G4").Value = 93094
$ws.Range("H4").Value = 80706
$ws.Range("B5").Value = 151291
$ws.Range("C5").Value = 216317
$ws.Range("D5").Value = 253374
$ws.Range("E5").Value = 424593
$ws.Range("F5").Value = 632222
$ws.Range("G5").Value = 131976
$ws.Range("H5").Value = 114794
$ws.Range("B6").Value = 260138
$ws.Range("C6").Value = 346515
$ws.Range("D6").Value = 456071
$ws.Range("E6").Value = 635678
$ws.Range("F6").Value = 866259
$ws.Range("G6").Value = 227939
$ws.Range("H6").Value = 198934
$ws.Range("B7").Value = 68031
$ws.Range("C7").Value = 93149
$ws.Range("D7").Value = 121748
$ws.Range("E7").Value = 230116
$ws.Range("F7").Value = 335731
$ws.Range("G7").Value = 64670
$ws.Range("H7").Value = 58202
$ws.Range("B8").Value = 69307
$ws.Range("C8").Value = 94069
$ws.Range("D8").Value = 121710
$ws.Range("E8").Value = 232012
$ws.Range("F8").Value = 341963
$ws.Range("G8").Value = 64361
$ws.Range("H8").Value = 60070
$ws.Range("B9").Value = 79546
$ws.Range("C9").Value = 106484
$ws.Range("D9").Value = 140522
$ws.Range("E9").Value = 269192
$ws.Range("F9").Value = 394851
$ws.Range("G9").Value = 75246
$ws.Range("H9").Value = 70482
$ws.Range("B10").Value = 96947
$ws.Range("C10").Value = 136282
$ws.Range("D10").Value = 167375
$ws.Range("E10").Value = 324089
$ws.Range("F10").Value = 462042
$ws.Range("G10").Value = 90398
$ws.Range("H10").Value = 84675
$ws.Range("B11").Value = 80636
$ws.Range("C11").Value = 103570
$ws.Range("D11").Value = 136303
$ws.Range("E11").Value = 231075
$ws.Range("F11").Value = 398895
$ws.Range("G11").Value = 70232
$ws.Range("H11").Value = 60579

# Update the active selection to D14 as in the diff
$ws.Range("D14").Select()
